$d = $word.ActiveDocument

# Array of old/new text pairs, applied in document order via Find/Replace All.
# MatchWholeWord is not used because these are standalone paragraph texts anyway,
# and every "old" value is unique in the document, so a simple ReplaceAll is safe.
$pairs = @(
    @("2026-02-11 Wednesday", "2026-02-12 Thursday"),
    @("92-24=", "28+24="),
    @("50-15=", "58+24="),
    @("60-46=", "29+63="),
    @("53-26=", "37+29="),
    @("92-67=", "29+5="),
    @("29+67=", "96-89="),
    @("9+49=", "80-67="),
    @("47+16=", "78-29="),
    @("93-37=", "60-38="),
    @("27+38=", "37+57="),
    @("48+44=", "50-31="),
    @("36-27=", "66-9="),
    @("40-14=", "85-58="),
    @("25+66=", "2+19="),
    @("43-37=", "56+19="),
    @("40-12=", "63-35="),
    @("70-66=", "22+19="),
    @("68+27=", "5+66="),
    @("86-37=", "49+39="),
    @("19+59=", "68+26="),
    @("33+49=", "61-54="),
    @("96-57=", "30-14="),
    @("6+27=", "90-53="),
    @("40-4=", "47+49="),
    @("36+7=", "91-33="),
    @("90-19=", "95-29="),
    @("13+79=", "90-3="),
    @("80-12=", "7+17="),
    @("59+38=", "30-5="),
    @("9+32=", "22+39="),
    @("70-11=", "51-24="),
    @("23+69=", "28+15="),
    @("90-36=", "4+69="),
    @("77+6=", "69+24="),
    @("18+74=", "90-2="),
    @("39+44=", "96-87="),
    @("57-18=", "39+57="),
    @("98-19=", "47+5="),
    @("8+3=", "18+14="),
    @("24+48=", "25+18="),
    @("69+14=", "58-29="),
    @("92-35=", "43-16="),
    @("39+19=", "72-39="),
    @("8+89=", "17+39="),
    @("65+28=", "82-29="),
    @("55+8=", "61-12="),
    @("26+56=", "96-19="),
    @("28+57=", "4+79="),
    @("81-25=", "36+57="),
    @("18+9=", "45+18="),
    @("45+38=", "50-5="),
    @("17+47=", "36+19="),
    @("19+62=", "19+43="),
    @("33+39=", "17+28="),
    @("33-28=", "93-7="),
    @("79+5=", "61-43="),
    @("6+6=", "74-16="),
    @("92-14=", "35+46="),
    @("81-66=", "90-73="),
    @("49+33=", "19+54="),
    @("71-44=", "17+36="),
    @("86-58=", "39+32="),
    @("71-56=", "13+58="),
    @("42-25=", "76-58="),
    @("51-43=", "70-53="),
    @("48+14=", "40-32="),
    @("9+72=", "91-38="),
    @("71-55=", "54-5="),
    @("87+6=", "35-8="),
    @("52-35=", "6+29="),
    @("29+58=", "45+29="),
    @("49+48=", "91-63="),
    @("79+3=", "58-29="),
    @("23-4=", "41-16="),
    @("85-46=", "42-35="),
    @("83-29=", "27+37="),
    @("48+27=", "57-9="),
    @("27+35=", "17+6="),
    @("88+9=", "92-49="),
    @("39+6=", "70-5="),
    @("82-69=", "38+14="),
    @("26-19=", "50-7="),
    @("39+4=", "7+69="),
    @("63-59=", "73+18="),
    @("76-49=", "83-75="),
    @("7+27=", "83-56="),
    @("19+65=", "57+15="),
    @("83-78=", "85-36="),
    @("19+75=", "39+14="),
    @("16+46=", "23+18="),
    @("7+59=", "80-54="),
    @("74-49=", "72+9="),
    @("55+17=", "59+16="),
    @("68-29=", "90-84="),
    @("39+59=", "9+74="),
    @("31-16=", "18+5="),
    @("71-46=", "70-57="),
    @("67-8=", "58+4="),
    @("53+19=", "35+8="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $result = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "WARNING: replacement not found for: $old"
    }
}

Write-Output "Done."
